# Auto update Excel log
# Appends new sensor log rows to the PIR, Humidity, Temperature, Proximity and mmWave sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('PIR')
$ws.Range('A46:F46').NumberFormat = '@'
$ws.Cells.Item(46, 1).Value = '2026-01-30'
$ws.Cells.Item(46, 2).Value = '12:55:28'
$ws.Cells.Item(46, 3).Value = '12:00'
$ws.Cells.Item(46, 4).Value = 'Bathroom'
$ws.Cells.Item(46, 5).Value = 'No Motion'
$ws.Cells.Item(46, 6).Value = 'Inactive'
$ws.Range('A47:F47').NumberFormat = '@'
$ws.Cells.Item(47, 1).Value = '2026-01-30'
$ws.Cells.Item(47, 2).Value = '12:55:30'
$ws.Cells.Item(47, 3).Value = '12:00'
$ws.Cells.Item(47, 4).Value = 'Bathroom'
$ws.Cells.Item(47, 5).Value = 'No Motion'
$ws.Cells.Item(47, 6).Value = 'Inactive'
$ws.Range('A48:F48').NumberFormat = '@'
$ws.Cells.Item(48, 1).Value = '2026-01-30'
$ws.Cells.Item(48, 2).Value = '12:55:35'
$ws.Cells.Item(48, 3).Value = '12:00'
$ws.Cells.Item(48, 4).Value = 'Bathroom'
$ws.Cells.Item(48, 5).Value = 'No Motion'
$ws.Cells.Item(48, 6).Value = 'Inactive'
$ws.Range('A49:F49').NumberFormat = '@'
$ws.Cells.Item(49, 1).Value = '2026-01-30'
$ws.Cells.Item(49, 2).Value = '12:55:40'
$ws.Cells.Item(49, 3).Value = '12:00'
$ws.Cells.Item(49, 4).Value = 'Bathroom'
$ws.Cells.Item(49, 5).Value = 'No Motion'
$ws.Cells.Item(49, 6).Value = 'Inactive'
$ws.Range('A50:F50').NumberFormat = '@'
$ws.Cells.Item(50, 1).Value = '2026-01-30'
$ws.Cells.Item(50, 2).Value = '12:55:45'
$ws.Cells.Item(50, 3).Value = '12:00'
$ws.Cells.Item(50, 4).Value = 'Bathroom'
$ws.Cells.Item(50, 5).Value = 'No Motion'
$ws.Cells.Item(50, 6).Value = 'Inactive'
$ws.Range('A51:F51').NumberFormat = '@'
$ws.Cells.Item(51, 1).Value = '2026-01-30'
$ws.Cells.Item(51, 2).Value = '12:55:50'
$ws.Cells.Item(51, 3).Value = '12:00'
$ws.Cells.Item(51, 4).Value = 'Bathroom'
$ws.Cells.Item(51, 5).Value = 'No Motion'
$ws.Cells.Item(51, 6).Value = 'Inactive'
$ws.Range('A52:F52').NumberFormat = '@'
$ws.Cells.Item(52, 1).Value = '2026-01-30'
$ws.Cells.Item(52, 2).Value = '12:55:53'
$ws.Cells.Item(52, 3).Value = '12:00'
$ws.Cells.Item(52, 4).Value = 'Living Room'
$ws.Cells.Item(52, 5).Value = 'RECOVERY_DETECTION'
$ws.Cells.Item(52, 6).Value = 'Inactive'
$ws.Range('A53:F53').NumberFormat = '@'
$ws.Cells.Item(53, 1).Value = '2026-01-30'
$ws.Cells.Item(53, 2).Value = '12:55:55'
$ws.Cells.Item(53, 3).Value = '12:00'
$ws.Cells.Item(53, 4).Value = 'Bathroom'
$ws.Cells.Item(53, 5).Value = 'No Motion'
$ws.Cells.Item(53, 6).Value = 'Inactive'
$ws.Range('A54:F54').NumberFormat = '@'
$ws.Cells.Item(54, 1).Value = '2026-01-30'
$ws.Cells.Item(54, 2).Value = '12:56:00'
$ws.Cells.Item(54, 3).Value = '12:00'
$ws.Cells.Item(54, 4).Value = 'Bathroom'
$ws.Cells.Item(54, 5).Value = 'No Motion'
$ws.Cells.Item(54, 6).Value = 'Inactive'
$ws.Range('A55:F55').NumberFormat = '@'
$ws.Cells.Item(55, 1).Value = '2026-01-30'
$ws.Cells.Item(55, 2).Value = '12:56:05'
$ws.Cells.Item(55, 3).Value = '12:00'
$ws.Cells.Item(55, 4).Value = 'Bathroom'
$ws.Cells.Item(55, 5).Value = 'No Motion'
$ws.Cells.Item(55, 6).Value = 'Inactive'
$ws.Range('A56:F56').NumberFormat = '@'
$ws.Cells.Item(56, 1).Value = '2026-01-30'
$ws.Cells.Item(56, 2).Value = '12:56:10'
$ws.Cells.Item(56, 3).Value = '12:00'
$ws.Cells.Item(56, 4).Value = 'Bathroom'
$ws.Cells.Item(56, 5).Value = 'No Motion'
$ws.Cells.Item(56, 6).Value = 'Inactive'
$ws.Range('A57:F57').NumberFormat = '@'
$ws.Cells.Item(57, 1).Value = '2026-01-30'
$ws.Cells.Item(57, 2).Value = '12:56:15'
$ws.Cells.Item(57, 3).Value = '12:00'
$ws.Cells.Item(57, 4).Value = 'Bathroom'
$ws.Cells.Item(57, 5).Value = 'No Motion'
$ws.Cells.Item(57, 6).Value = 'Inactive'
$ws.Range('A58:F58').NumberFormat = '@'
$ws.Cells.Item(58, 1).Value = '2026-01-30'
$ws.Cells.Item(58, 2).Value = '12:56:20'
$ws.Cells.Item(58, 3).Value = '12:00'
$ws.Cells.Item(58, 4).Value = 'Bathroom'
$ws.Cells.Item(58, 5).Value = 'No Motion'
$ws.Cells.Item(58, 6).Value = 'Inactive'
$ws.Range('A59:F59').NumberFormat = '@'
$ws.Cells.Item(59, 1).Value = '2026-01-30'
$ws.Cells.Item(59, 2).Value = '12:56:25'
$ws.Cells.Item(59, 3).Value = '12:00'
$ws.Cells.Item(59, 4).Value = 'Bathroom'
$ws.Cells.Item(59, 5).Value = 'No Motion'
$ws.Cells.Item(59, 6).Value = 'Inactive'

$ws = $wb.Worksheets.Item('Humidity')
$ws.Range('A40:F40').NumberFormat = '@'
$ws.Cells.Item(40, 1).Value = '2026-01-30'
$ws.Cells.Item(40, 2).Value = '12:55:28'
$ws.Cells.Item(40, 3).Value = '12:00'
$ws.Cells.Item(40, 4).Value = 'Bathroom'
$ws.Cells.Item(40, 5).Value = '87.6%'
$ws.Cells.Item(40, 6).Value = 'Active'
$ws.Range('A41:F41').NumberFormat = '@'
$ws.Cells.Item(41, 1).Value = '2026-01-30'
$ws.Cells.Item(41, 2).Value = '12:55:29'
$ws.Cells.Item(41, 3).Value = '12:00'
$ws.Cells.Item(41, 4).Value = 'Bathroom'
$ws.Cells.Item(41, 5).Value = '87.6%'
$ws.Cells.Item(41, 6).Value = 'Active'
$ws.Range('A42:F42').NumberFormat = '@'
$ws.Cells.Item(42, 1).Value = '2026-01-30'
$ws.Cells.Item(42, 2).Value = '12:55:30'
$ws.Cells.Item(42, 3).Value = '12:00'
$ws.Cells.Item(42, 4).Value = 'Bathroom'
$ws.Cells.Item(42, 5).Value = '86.6%'
$ws.Cells.Item(42, 6).Value = 'Active'
$ws.Range('A43:F43').NumberFormat = '@'
$ws.Cells.Item(43, 1).Value = '2026-01-30'
$ws.Cells.Item(43, 2).Value = '12:55:34'
$ws.Cells.Item(43, 3).Value = '12:00'
$ws.Cells.Item(43, 4).Value = 'Bathroom'
$ws.Cells.Item(43, 5).Value = '86.4%'
$ws.Cells.Item(43, 6).Value = 'Active'
$ws.Range('A44:F44').NumberFormat = '@'
$ws.Cells.Item(44, 1).Value = '2026-01-30'
$ws.Cells.Item(44, 2).Value = '12:55:38'
$ws.Cells.Item(44, 3).Value = '12:00'
$ws.Cells.Item(44, 4).Value = 'Bathroom'
$ws.Cells.Item(44, 5).Value = '86.7%'
$ws.Cells.Item(44, 6).Value = 'Active'
$ws.Range('A45:F45').NumberFormat = '@'
$ws.Cells.Item(45, 1).Value = '2026-01-30'
$ws.Cells.Item(45, 2).Value = '12:55:42'
$ws.Cells.Item(45, 3).Value = '12:00'
$ws.Cells.Item(45, 4).Value = 'Bathroom'
$ws.Cells.Item(45, 5).Value = '86.1%'
$ws.Cells.Item(45, 6).Value = 'Active'
$ws.Range('A46:F46').NumberFormat = '@'
$ws.Cells.Item(46, 1).Value = '2026-01-30'
$ws.Cells.Item(46, 2).Value = '12:55:46'
$ws.Cells.Item(46, 3).Value = '12:00'
$ws.Cells.Item(46, 4).Value = 'Bathroom'
$ws.Cells.Item(46, 5).Value = '87.6%'
$ws.Cells.Item(46, 6).Value = 'Active'
$ws.Range('A47:F47').NumberFormat = '@'
$ws.Cells.Item(47, 1).Value = '2026-01-30'
$ws.Cells.Item(47, 2).Value = '12:55:50'
$ws.Cells.Item(47, 3).Value = '12:00'
$ws.Cells.Item(47, 4).Value = 'Bathroom'
$ws.Cells.Item(47, 5).Value = '86.7%'
$ws.Cells.Item(47, 6).Value = 'Active'
$ws.Range('A48:F48').NumberFormat = '@'
$ws.Cells.Item(48, 1).Value = '2026-01-30'
$ws.Cells.Item(48, 2).Value = '12:55:54'
$ws.Cells.Item(48, 3).Value = '12:00'
$ws.Cells.Item(48, 4).Value = 'Bathroom'
$ws.Cells.Item(48, 5).Value = '87.7%'
$ws.Cells.Item(48, 6).Value = 'Active'
$ws.Range('A49:F49').NumberFormat = '@'
$ws.Cells.Item(49, 1).Value = '2026-01-30'
$ws.Cells.Item(49, 2).Value = '12:56:06'
$ws.Cells.Item(49, 3).Value = '12:00'
$ws.Cells.Item(49, 4).Value = 'Bathroom'
$ws.Cells.Item(49, 5).Value = '87.7%'
$ws.Cells.Item(49, 6).Value = 'Active'
$ws.Range('A50:F50').NumberFormat = '@'
$ws.Cells.Item(50, 1).Value = '2026-01-30'
$ws.Cells.Item(50, 2).Value = '12:56:10'
$ws.Cells.Item(50, 3).Value = '12:00'
$ws.Cells.Item(50, 4).Value = 'Bathroom'
$ws.Cells.Item(50, 5).Value = '86.7%'
$ws.Cells.Item(50, 6).Value = 'Active'
$ws.Range('A51:F51').NumberFormat = '@'
$ws.Cells.Item(51, 1).Value = '2026-01-30'
$ws.Cells.Item(51, 2).Value = '12:56:18'
$ws.Cells.Item(51, 3).Value = '12:00'
$ws.Cells.Item(51, 4).Value = 'Bathroom'
$ws.Cells.Item(51, 5).Value = '87.7%'
$ws.Cells.Item(51, 6).Value = 'Active'
$ws.Range('A52:F52').NumberFormat = '@'
$ws.Cells.Item(52, 1).Value = '2026-01-30'
$ws.Cells.Item(52, 2).Value = '12:56:27'
$ws.Cells.Item(52, 3).Value = '12:00'
$ws.Cells.Item(52, 4).Value = 'Bathroom'
$ws.Cells.Item(52, 5).Value = '87.7%'
$ws.Cells.Item(52, 6).Value = 'Active'

$ws = $wb.Worksheets.Item('Temperature')
$ws.Range('A40:F40').NumberFormat = '@'
$ws.Cells.Item(40, 1).Value = '2026-01-30'
$ws.Cells.Item(40, 2).Value = '12:55:28'
$ws.Cells.Item(40, 3).Value = '12:00'
$ws.Cells.Item(40, 4).Value = 'Bathroom'
$ws.Cells.Item(40, 5).Value = '22.7C'
$ws.Cells.Item(40, 6).Value = 'Active'
$ws.Range('A41:F41').NumberFormat = '@'
$ws.Cells.Item(41, 1).Value = '2026-01-30'
$ws.Cells.Item(41, 2).Value = '12:55:29'
$ws.Cells.Item(41, 3).Value = '12:00'
$ws.Cells.Item(41, 4).Value = 'Bathroom'
$ws.Cells.Item(41, 5).Value = '22.7C'
$ws.Cells.Item(41, 6).Value = 'Active'
$ws.Range('A42:F42').NumberFormat = '@'
$ws.Cells.Item(42, 1).Value = '2026-01-30'
$ws.Cells.Item(42, 2).Value = '12:55:30'
$ws.Cells.Item(42, 3).Value = '12:00'
$ws.Cells.Item(42, 4).Value = 'Bathroom'
$ws.Cells.Item(42, 5).Value = '22.6C'
$ws.Cells.Item(42, 6).Value = 'Active'
$ws.Range('A43:F43').NumberFormat = '@'
$ws.Cells.Item(43, 1).Value = '2026-01-30'
$ws.Cells.Item(43, 2).Value = '12:55:35'
$ws.Cells.Item(43, 3).Value = '12:00'
$ws.Cells.Item(43, 4).Value = 'Bathroom'
$ws.Cells.Item(43, 5).Value = '22.6C'
$ws.Cells.Item(43, 6).Value = 'Active'
$ws.Range('A44:F44').NumberFormat = '@'
$ws.Cells.Item(44, 1).Value = '2026-01-30'
$ws.Cells.Item(44, 2).Value = '12:55:38'
$ws.Cells.Item(44, 3).Value = '12:00'
$ws.Cells.Item(44, 4).Value = 'Bathroom'
$ws.Cells.Item(44, 5).Value = '22.7C'
$ws.Cells.Item(44, 6).Value = 'Active'
$ws.Range('A45:F45').NumberFormat = '@'
$ws.Cells.Item(45, 1).Value = '2026-01-30'
$ws.Cells.Item(45, 2).Value = '12:55:43'
$ws.Cells.Item(45, 3).Value = '12:00'
$ws.Cells.Item(45, 4).Value = 'Bathroom'
$ws.Cells.Item(45, 5).Value = '22.6C'
$ws.Cells.Item(45, 6).Value = 'Active'
$ws.Range('A46:F46').NumberFormat = '@'
$ws.Cells.Item(46, 1).Value = '2026-01-30'
$ws.Cells.Item(46, 2).Value = '12:55:47'
$ws.Cells.Item(46, 3).Value = '12:00'
$ws.Cells.Item(46, 4).Value = 'Bathroom'
$ws.Cells.Item(46, 5).Value = '22.6C'
$ws.Cells.Item(46, 6).Value = 'Active'
$ws.Range('A47:F47').NumberFormat = '@'
$ws.Cells.Item(47, 1).Value = '2026-01-30'
$ws.Cells.Item(47, 2).Value = '12:55:51'
$ws.Cells.Item(47, 3).Value = '12:00'
$ws.Cells.Item(47, 4).Value = 'Bathroom'
$ws.Cells.Item(47, 5).Value = '22.6C'
$ws.Cells.Item(47, 6).Value = 'Active'
$ws.Range('A48:F48').NumberFormat = '@'
$ws.Cells.Item(48, 1).Value = '2026-01-30'
$ws.Cells.Item(48, 2).Value = '12:55:55'
$ws.Cells.Item(48, 3).Value = '12:00'
$ws.Cells.Item(48, 4).Value = 'Bathroom'
$ws.Cells.Item(48, 5).Value = '22.6C'
$ws.Cells.Item(48, 6).Value = 'Active'
$ws.Range('A49:F49').NumberFormat = '@'
$ws.Cells.Item(49, 1).Value = '2026-01-30'
$ws.Cells.Item(49, 2).Value = '12:56:07'
$ws.Cells.Item(49, 3).Value = '12:00'
$ws.Cells.Item(49, 4).Value = 'Bathroom'
$ws.Cells.Item(49, 5).Value = '22.6C'
$ws.Cells.Item(49, 6).Value = 'Active'
$ws.Range('A50:F50').NumberFormat = '@'
$ws.Cells.Item(50, 1).Value = '2026-01-30'
$ws.Cells.Item(50, 2).Value = '12:56:11'
$ws.Cells.Item(50, 3).Value = '12:00'
$ws.Cells.Item(50, 4).Value = 'Bathroom'
$ws.Cells.Item(50, 5).Value = '22.6C'
$ws.Cells.Item(50, 6).Value = 'Active'
$ws.Range('A51:F51').NumberFormat = '@'
$ws.Cells.Item(51, 1).Value = '2026-01-30'
$ws.Cells.Item(51, 2).Value = '12:56:19'
$ws.Cells.Item(51, 3).Value = '12:00'
$ws.Cells.Item(51, 4).Value = 'Bathroom'
$ws.Cells.Item(51, 5).Value = '22.6C'
$ws.Cells.Item(51, 6).Value = 'Active'
$ws.Range('A52:F52').NumberFormat = '@'
$ws.Cells.Item(52, 1).Value = '2026-01-30'
$ws.Cells.Item(52, 2).Value = '12:56:27'
$ws.Cells.Item(52, 3).Value = '12:00'
$ws.Cells.Item(52, 4).Value = 'Bathroom'
$ws.Cells.Item(52, 5).Value = '22.6C'
$ws.Cells.Item(52, 6).Value = 'Active'

$ws = $wb.Worksheets.Item('Proximity')
$ws.Range('A18:F18').NumberFormat = '@'
$ws.Cells.Item(18, 1).Value = '2026-01-30'
$ws.Cells.Item(18, 2).Value = '12:55:30'
$ws.Cells.Item(18, 3).Value = '12:00'
$ws.Cells.Item(18, 4).Value = 'Bathroom Door'
$ws.Cells.Item(18, 5).Value = 'ENTER'
$ws.Cells.Item(18, 6).Value = 'User ENTERED Bathroom'
$ws.Range('A19:F19').NumberFormat = '@'
$ws.Cells.Item(19, 1).Value = '2026-01-30'
$ws.Cells.Item(19, 2).Value = '12:55:45'
$ws.Cells.Item(19, 3).Value = '12:00'
$ws.Cells.Item(19, 4).Value = 'Bathroom Door'
$ws.Cells.Item(19, 5).Value = 'EXIT'
$ws.Cells.Item(19, 6).Value = 'User EXITED Bathroom'
$ws.Range('A20:F20').NumberFormat = '@'
$ws.Cells.Item(20, 1).Value = '2026-01-30'
$ws.Cells.Item(20, 2).Value = '12:55:51'
$ws.Cells.Item(20, 3).Value = '12:00'
$ws.Cells.Item(20, 4).Value = 'Bathroom Door'
$ws.Cells.Item(20, 5).Value = 'ENTER'
$ws.Cells.Item(20, 6).Value = 'User ENTERED Bathroom'
$ws.Range('A21:F21').NumberFormat = '@'
$ws.Cells.Item(21, 1).Value = '2026-01-30'
$ws.Cells.Item(21, 2).Value = '12:55:57'
$ws.Cells.Item(21, 3).Value = '12:00'
$ws.Cells.Item(21, 4).Value = 'Bathroom Door'
$ws.Cells.Item(21, 5).Value = 'EXIT'
$ws.Cells.Item(21, 6).Value = 'User EXITED Bathroom'
$ws.Range('A22:F22').NumberFormat = '@'
$ws.Cells.Item(22, 1).Value = '2026-01-30'
$ws.Cells.Item(22, 2).Value = '12:56:11'
$ws.Cells.Item(22, 3).Value = '12:00'
$ws.Cells.Item(22, 4).Value = 'Bathroom Door'
$ws.Cells.Item(22, 5).Value = 'ENTER'
$ws.Cells.Item(22, 6).Value = 'User ENTERED Bathroom'
$ws.Range('A23:F23').NumberFormat = '@'
$ws.Cells.Item(23, 1).Value = '2026-01-30'
$ws.Cells.Item(23, 2).Value = '12:56:18'
$ws.Cells.Item(23, 3).Value = '12:00'
$ws.Cells.Item(23, 4).Value = 'Bathroom Door'
$ws.Cells.Item(23, 5).Value = 'EXIT'
$ws.Cells.Item(23, 6).Value = 'User EXITED Bathroom'
$ws.Range('A24:F24').NumberFormat = '@'
$ws.Cells.Item(24, 1).Value = '2026-01-30'
$ws.Cells.Item(24, 2).Value = '12:56:27'
$ws.Cells.Item(24, 3).Value = '12:00'
$ws.Cells.Item(24, 4).Value = 'Bathroom Door'
$ws.Cells.Item(24, 5).Value = 'ENTER'
$ws.Cells.Item(24, 6).Value = 'User ENTERED Bathroom'

$ws = $wb.Worksheets.Item('mmWave')
$ws.Range('A21:F21').NumberFormat = '@'
$ws.Cells.Item(21, 1).Value = '2026-01-30'
$ws.Cells.Item(21, 2).Value = '12:55:53'
$ws.Cells.Item(21, 3).Value = '12:00'
$ws.Cells.Item(21, 4).Value = 'Living Room'
$ws.Cells.Item(21, 5).Value = 'PRESENCE_DETECTED'
$ws.Cells.Item(21, 6).Value = 'Active'
$ws.Range('A22:F22').NumberFormat = '@'
$ws.Cells.Item(22, 1).Value = '2026-01-30'
$ws.Cells.Item(22, 2).Value = '12:56:14'
$ws.Cells.Item(22, 3).Value = '12:00'
$ws.Cells.Item(22, 4).Value = 'Living Room'
$ws.Cells.Item(22, 5).Value = 'FALL_DETECTED'
$ws.Cells.Item(22, 6).Value = 'EMERGENCY'

